# Insert a new data row at row 252 (pushing existing rows 252..353 down to 253..354)
# and populate it with the new "Ají" record (weekly price update commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(252).Insert()

$ws.Range("A252").Value = 4
$ws.Range("B252").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C252").Value = "Los Lagos"
$ws.Range("D252").Value = 44917
$ws.Range("E252").Value = 10
$ws.Range("F252").Value = 100112021
$ws.Range("G252").Value = "Ají"
$ws.Range("H252").Value = "Inferno"
$ws.Range("I252").Value = "Primera"
$ws.Range("J252").Value = 70
$ws.Range("K252").Value = 22000
$ws.Range("L252").Value = 22000
$ws.Range("M252").Value = 22000
$ws.Range("N252").Value = "`$/caja 10 kilos"
$ws.Range("O252").Value = "Región de Arica y Parinacota"
$ws.Range("P252").Value = 2200
$ws.Range("Q252").Value = 10
$ws.Range("R252").Value = "Hortaliza"
